$d = $word.ActiveDocument

# Step 1: append a temp paragraph at the very end of the document containing our new run's text.
# We build it with no explicit character formatting so it has no <w:rPr> residue.
$endRng = $d.Content
$endRng.Collapse(0)
$markBeforeTemp = $endRng.Start
$endRng.InsertParagraphAfter()
$endRng.Collapse(0)
$tempStart = $endRng.Start
$endRng.InsertAfter(" How many more years of data do we need to achieve a reasonable type II error.")
$tempEnd = $endRng.End
$tempParaEnd = $tempEnd + 1   # include the trailing paragraph mark
Write-Host "temp: markBeforeTemp=$markBeforeTemp tempStart=$tempStart tempEnd=$tempEnd tempParaEnd=$tempParaEnd"

# Step 2: locate insertion point in the target paragraph
$rng = $d.Content
$found = $rng.Find.Execute("How severe do we expect it to be?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Find result: $found"
$rng.Collapse(0)
$insStart = $rng.Start
Write-Host "insStart = $insStart"

# Step 3: copy the temp run's FormattedText (preserves it as its own run, no rPr) to the target location
$destRng = $d.Range($insStart, $insStart)
$srcRng = $d.Range($tempStart, $tempEnd)
$destRng.FormattedText = $srcRng
Write-Host "Done assigning"

# Step 4: remove the temp paragraph we appended at the end (text + paragraph mark).
# Everything from markBeforeTemp onward shifted forward by the length of the text we just inserted.
$shift = $tempEnd - $tempStart
$delStart = $markBeforeTemp + $shift
$delEnd = $tempParaEnd + $shift
$delRng = $d.Range($delStart, $delEnd)
Write-Host "delRng=[$($delRng.Text)]"
$delRng.Delete()
Write-Host "Done cleanup. Doc end=$($d.Content.End)  Paragraphs.Count=$($d.Paragraphs.Count)"

# Step 5: merge the "W" / "hat sample size..." runs into a single run by replacing the
# full sentence with itself (Find/Replace naturally collapses it into one run).
$rng2 = $d.Content
$sentence = "What sample size is needed for RS of age 0 fish to be a reliable predictor of TLF. This has significant management implications, for example the relationship between RS and TLF determines the efficacy of labor at screw traps."
$found2 = $rng2.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)
Write-Host "Find2 result: $found2"
